# Update "想去人数" (interested-count) figures on the "展览" and "全部类型"
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 3397
    $ws.Range("F5").Value = 1585
    $ws.Range("F6").Value = 60
}
